$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.325.11'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '3.345.98'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'249.02"
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('D6').Value = "'652.24"
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('D7').Value = "'1.39"
$ws.Range('E7').Value = '  -9.45%  '
$ws.Range('D8').Value = "'0.417"
$ws.Range('E8').Value = '  -11.91%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -5.67%  '
$ws.Range('D11').Value = '3.346.62'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('E12').Value = '  -5.22%  '
$ws.Range('D13').Value = "'40.21"
$ws.Range('E13').Value = '  -5.74%  '
$ws.Range('D14').Value = '97.051.10'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = "'6.04"
$ws.Range('E15').Value = '  +4.94%  '
$ws.Range('D16').Value = "'0.0000254"
$ws.Range('E16').Value = '  -7.16%  '
$ws.Range('D17').Value = '3.968.75'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = "'8.54"
$ws.Range('E18').Value = '  +5.77%  '
$ws.Range('D19').Value = '3.348.30'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').Value = "'0.544"
$ws.Range('E20').Value = '  +23.74%  '
$ws.Range('D21').Value = "'16.78"
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('D22').Value = "'10.71"
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').Value = "'497.55"
$ws.Range('E23').Value = '  -6.13%  '
$ws.Range('D24').Value = "'3.33"
$ws.Range('E24').Value = '  -6.28%  '
$ws.Range('D25').Value = "'0.0000200"
$ws.Range('E25').Value = '  -7.71%  '
$ws.Range('D26').Value = "'6.31"
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = "'95.89"
$ws.Range('E27').Value = '  -6.98%  '
$ws.Range('D28').Value = "'12.02"
$ws.Range('E28').Value = '  -5.17%  '
$ws.Range('D29').Value = '3.529.64'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').Value = "'0.147"
$ws.Range('E30').Value = '  -4.21%  '
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').Value = "'11.01"
$ws.Range('E32').Value = '  -4.07%  '
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').Value = "'2.47"
$ws.Range('E34').Value = '  +15.48%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = "'0.549"
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('D37').Value = "'28.50"
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('D39').Value = "'1.44"
$ws.Range('E39').Value = '  +7.62%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = "'508.65"
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('E42').Value = '  -6.32%  '
$ws.Range('D43').Value = "'24.61"
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('D44').Value = "'0.838"
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  +8.68%  '
$ws.Range('E46').Value = '  -5.24%  '
$ws.Range('D47').Value = "'0.0412"
$ws.Range('E47').Value = '  -5.30%  '
$ws.Range('D48').Value = "'5.54"
$ws.Range('E48').Value = '  +5.93%  '
$ws.Range('D49').Value = "'1.64"
$ws.Range('E49').Value = '  +5.93%  '
$ws.Range('D50').Value = "'53.05"
$ws.Range('E50').Value = '  +5.01%  '
$ws.Range('E51').Value = '  -9.25%  '
